# Proyecto trimestral: proceso de facturacion completado.
# Adds three new product rows (Cocacola, Colombiana, Margaritas) to the
# "Inventario" sheet, below the existing products.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=cod_articulo, B=nombre_producto, C=disponibilidad,
#          D=valor_unitario, E=iva_producto, F=categoria, G=vencimiento_producto
# Force text format so codes like "0004" and dates like "2024-09-11" are
# stored as text (matching the rest of the table) instead of being
# auto-converted to numbers/dates.
$ws.Range("A5:G7").NumberFormat = "@"

$ws.Range("A5").Value = "0004"
$ws.Range("B5").Value = "Cocacola"
$ws.Range("C5").Value = "500"
$ws.Range("D5").Value = "3500"
$ws.Range("E5").Value = "19"
$ws.Range("F5").Value = "2"
$ws.Range("G5").Value = "2024-09-11"

$ws.Range("A6").Value = "0005"
$ws.Range("B6").Value = "Colombiana"
$ws.Range("C6").Value = "300"
$ws.Range("D6").Value = "2800"
$ws.Range("E6").Value = "19"
$ws.Range("F6").Value = "2"
$ws.Range("G6").Value = "2024-10-25"

$ws.Range("A7").Value = "0006"
$ws.Range("B7").Value = "Margaritas"
$ws.Range("C7").Value = "2000"
$ws.Range("D7").Value = "400"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "1"
$ws.Range("G7").Value = "2021-11-20"

# Match the saved selection state from the edit (cell below the new rows).
$ws.Range("E8").Select()
